$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting N,O,P -> O,P,Q
$ws.Columns("N:N").Insert()

# The newly inserted column adopts a width close to the former column M
$ws.Columns("N:N").ColumnWidth = 9.83

# Make "Repayment schedule" the active/selected sheet (it moves from Summary)
$ws.Activate()
$ws.Range("S5").Select() | Out-Null
